$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.612.62'
$ws.Range("E2").Value = '  +6.77%  '
$ws.Range("D3").Value = '3.508.50'
$ws.Range("E3").Value = '  +9.38%  '
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").Value = '''190.87'
$ws.Range("E5").Value = '  +10.49%  '
$ws.Range("D6").Value = '''556.18'
$ws.Range("E6").Value = '  +7.75%  '
$ws.Range("D7").Value = '3.496.64'
$ws.Range("E7").Value = '  +9.19%  '
$ws.Range("D8").Value = '''0.608'
$ws.Range("E8").Value = '  +3.39%  '
$ws.Range("E9").Value = '  -0.10%  '
$ws.Range("D10").Value = '''0.637'
$ws.Range("E10").Value = '  +6.90%  '
$ws.Range("D11").Value = '''56.20'
$ws.Range("E11").Value = '  +7.99%  '
$ws.Range("D12").Value = '''0.150'
$ws.Range("E12").Value = '  +15.98%  '
$ws.Range("D13").Value = '''0.0000269'
$ws.Range("E13").Value = '  +7.98%  '
$ws.Range("D14").Value = '''9.43'
$ws.Range("E14").Value = '  +6.73%  '
$ws.Range("D15").Value = '4.098.74'
$ws.Range("E15").Value = '  +8.92%  '
$ws.Range("D16").Value = '3.522.28'
$ws.Range("E16").Value = '  +8.66%  '
$ws.Range("D17").Value = '''0.121'
$ws.Range("E17").Value = '  +6.58%  '
$ws.Range("D18").Value = '''18.30'
$ws.Range("E18").Value = '  +7.46%  '
$ws.Range("D19").Value = '66.824.49'
$ws.Range("E19").Value = '  +6.70%  '
$ws.Range("D20").Value = '''11.88'
$ws.Range("E20").Value = '  +9.37%  '
$ws.Range("D21").Value = '''0.997'
$ws.Range("E21").Value = '  +4.81%  '
$ws.Range("D22").Value = '''407.53'
$ws.Range("E22").Value = '  +12.51%  '
$ws.Range("D23").Value = '''3.96'
$ws.Range("E23").Value = '  +7.24%  '
$ws.Range("D24").Value = '''85.35'
$ws.Range("E24").Value = '  +7.07%  '
$ws.Range("D25").Value = '''4.21'
$ws.Range("E25").Value = '  +8.95%  '
$ws.Range("D26").Value = '''11.18'
$ws.Range("E26").Value = '  +2.73%  '
$ws.Range("D27").Value = '''2.94'
$ws.Range("E27").Value = '  +13.51%  '
$ws.Range("E28").Value = '  +0.33%  '
$ws.Range("D29").Value = '''11.94'
$ws.Range("E29").Value = '  +7.83%  '
$ws.Range("D30").Value = '''8.85'
$ws.Range("E30").Value = '  +9.54%  '
$ws.Range("D31").Value = '''30.42'
$ws.Range("E31").Value = '  +8.65%  '
$ws.Range("D32").Value = '''665.83'
$ws.Range("E32").Value = '  +3.18%  '
$ws.Range("D33").Value = '''6.72'
$ws.Range("E33").Value = '  +6.25%  '
$ws.Range("D34").Value = '''11.78'
$ws.Range("E34").Value = '  +7.01%  '
$ws.Range("D35").Value = '''0.111'
$ws.Range("E35").Value = '  +7.82%  '
$ws.Range("D36").Value = '''59.59'
$ws.Range("E36").Value = '  +4.54%  '
$ws.Range("D37").Value = '''38.84'
$ws.Range("E37").Value = '  +7.62%  '
$ws.Range("D38").Value = '0.0₃0815'
$ws.Range("E38").Value = '  +16.25%  '
$ws.Range("E39").Value = '  -0.02%  '
$ws.Range("D40").Value = '''0.394'
$ws.Range("E40").Value = '  +6.45%  '
$ws.Range("D41").Value = '''0.137'
$ws.Range("E41").Value = '  +13.39%  '
$ws.Range("D42").Value = '''3.37'
$ws.Range("E42").Value = '  +22.18%  '
$ws.Range("E43").Value = '  +20.35%  '
$ws.Range("D44").Value = '''1.00'
$ws.Range("E44").Value = '  -0.14%  '
$ws.Range("D45").Value = '3.011.34'
$ws.Range("E45").Value = '  +5.07%  '
$ws.Range("D46").Value = '''2.65'
$ws.Range("E46").Value = '  +7.44%  '
$ws.Range("D47").Value = '''3.35'
$ws.Range("E47").Value = '  +13.15%  '
$ws.Range("D48").Value = '''0.0417'
$ws.Range("E48").Value = '  +8.42%  '
$ws.Range("D49").Value = '''9.08'
$ws.Range("E49").Value = '  +19.81%  '
$ws.Range("D50").Value = '''2.73'
$ws.Range("E50").Value = '  +4.41%  '
$ws.Range("D51").Value = '''0.130'
$ws.Range("E51").Value = '  +6.91%  '
